$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3: EFT - Şube
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("G3").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("H3").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("J3").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# Row 4: EFT - ATM
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("G4").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("H4").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("J4").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# Row 5: EFT - Mobil
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("G5").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("H5").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("J5").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# Row 6: DÜZENLİ EFT
$ws.Range("D6").Value = ""
$ws.Range("E6").Value = "6,09 TL - 12,19 TL - 152,35 TL"
$ws.Range("G6").Value = "6,09 TL - 12,19 TL - 152,35 TL"
$ws.Range("H6").Value = "6,09 TL - 12,19 TL - 152,35 TL"
$ws.Range("J6").Value = "6,09 TL - 12,19 TL - 152,35 TL"

# Row 8: HAVALE - Şube
$ws.Range("D8").Value = ""
$ws.Range("E8").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("G8").Value = "15,23 TL - 30,47 TL - 211,05 TL"
$ws.Range("H8").Value = "15,23 TL - 30,47 TL - 304,72 TL"
$ws.Range("J8").Value = "15,23 TL - 30,47 TL - 304,72 TL"

# Row 9: HAVALE - ATM
$ws.Range("D9").Value = ""
$ws.Range("E9").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("G9").Value = "15,23 TL - 30,47 TL - 211,05 TL"
$ws.Range("H9").Value = "15,23 TL - 30,47 TL - 304,72 TL"
$ws.Range("J9").Value = "15,23 TL - 30,47 TL - 304,72 TL"

# Row 10: HAVALE - Mobil
$ws.Range("D10").Value = ""
$ws.Range("E10").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("G10").Value = "15,23 TL - 30,47 TL - 211,05 TL"
$ws.Range("H10").Value = "15,23 TL - 30,47 TL - 304,72 TL"
$ws.Range("J10").Value = "15,23 TL - 30,47 TL - 304,72 TL"

# Row 11: DÜZENLİ HAVALE
$ws.Range("D11").Value = ""
$ws.Range("E11").Value = "3,04 TL - 6,09 TL - 76,17 TL"
$ws.Range("G11").Value = "3,04 TL - 6,09 TL - 76,17 TL"
$ws.Range("H11").Value = "3,05 TL - 6,1 TL - 76,18 TL"
$ws.Range("J11").Value = "3,05 TL - 6,1 TL - 76,18 TL"

# Row 12: GİDEN SWIFT
$ws.Range("D12").Value = ""
$ws.Range("G12").Value = "Şube (Kasadan): %0,5; Şube (Hesaptan): %0,75; İnternet: 15 USD"

# Row 13: GELEN SWIFT
$ws.Range("C13").Value = "Hesaba: Asgari 0 TL | Azami 0,94 TL"
$ws.Range("D13").Value = ""
$ws.Range("E13").Value = "Hesaba: Asgari 1 TL | Azami 851,5 TL"
$ws.Range("H13").Value = "Hesaba: Asgari 1 TL | Azami 6,09 TL"
$ws.Range("J13").Value = "Hesaba: Asgari 1 TL | Azami 995,5 TL"

# Row 14: GİDEN SWIFT - Mobil
$ws.Range("D14").Value = ""
$ws.Range("E14").Value = "1.660 TL - 1.660 TL"
$ws.Range("G14").Value = "6.300 TL - 6,09 TL"
$ws.Range("H14").Value = "2.100 TL - 4.300 TL"
$ws.Range("J14").Value = "1.188 TL - 593 TL"
